$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "West"
$ws.Range("B2").Value = "N"

$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B2").Select()
